# Update the "K" column (G) values for rows 2-7 in Sheet1.
# These values correspond to the strike-count ("Strike#") -> "K" rework
# described in the commit message: regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 0
